$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Locate data rows by country name in column A ---
$brasilCell = $ws.Range("A1:A1000").Find("Brasil")
$brasilRow = $brasilCell.Row

$jamaicaCell = $ws.Range("A1:A1000").Find("Jamaica")
$jamaicaRow = $jamaicaCell.Row

# --- Update Brasil stats (row stays in place; ranking unaffected) ---
$ws.Cells.Item($brasilRow, 2).Value = 63100   # Casos totales
$ws.Cells.Item($brasilRow, 3).Value = 241     # Nuevos casos
$ws.Cells.Item($brasilRow, 5).Value = 28662   # Recuperados
$ws.Cells.Item($brasilRow, 7).Value = 15      # Muertes hoy
$ws.Cells.Item($brasilRow, 8).Value = 4286    # Muertes

# --- Update Jamaica stats ---
$ws.Cells.Item($jamaicaRow, 2).Value = 350    # Casos totales
$ws.Cells.Item($jamaicaRow, 3).Value = 45     # Nuevos casos
$ws.Cells.Item($jamaicaRow, 4).Value = 28     # Casos activos
$ws.Cells.Item($jamaicaRow, 5).Value = 315    # Recuperados
$ws.Cells.Item($jamaicaRow, 8).Value = 7      # Muertes

# --- Re-sort the whole country table descending by "Casos totales" (col B) ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$dataRange = $ws.Range("A4:H" + $lastRow)
$sortKey = $ws.Range("B4:B" + $lastRow)

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 2, 0, 0)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()
